$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the diff. "Price" (column D) values are
# free-form text that must stay text even when they look numeric
# (e.g. 577.67, 1.60), so format the cell as Text before assigning,
# then restore the default (Normal) style so no stray formatting
# is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.430.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.503.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.22%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.86%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.522"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.502.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = "  -7.17%  "
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.962.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.377.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.503.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.65%  "
$ws.Range("E19").Value = "  -7.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.632.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0902"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "462.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.59%  "
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +3.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.09%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.318"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -14.67%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.70%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.65%  "
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("E49").Value = "  -4.19%  "
$ws.Range("E50").Value = "  -5.10%  "
$ws.Range("E51").Value = "  -1.87%  "
